$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 12: intersection formula -> #NULL!
$ws1.Range("K12").Formula = "=J8:K10 L6:M7"

# Row 13-14: shared string + text-concat formulas
$ws1.Range("C13").Value = "ca"
$ws1.Range("E13").Formula = '="ciao"&B13&" cc" &C13'
$ws1.Range("E14").Formula = '="ciao"&B13&" cc" &C13'

# Row 15: plain number, plus a quote-prefixed numeric cell (style only)
$ws1.Range("B15").Value = 1
$ws1.Range("C15").Value = 0

# Row 17: blank cell that will carry the quote-prefix style (built below)
$ws1.Range("C17").Value = "'x"
$ws1.Range("C17").Value = ""

# Copy C17's (quote-prefix) formatting onto C15 without touching its value
$ws1.Range("C17").Copy()
$ws1.Range("C15").PasteSpecial(-4122)

# Row 16-18: 3x3 legacy CSE array formula spilling G11:H12 (2x2) -> #N/A overflow
$ws1.Range("H16:J18").FormulaArray = "=G11:H12"

# Row 19: simple formula
$ws1.Range("E19").Formula = "=D16*3"

# Row 20-24: 5x1 legacy CSE array formula
$ws1.Range("H20:H24").FormulaArray = "=IF(G9:G12<>H10:H13,1,0)"

# Row 21: simple formula
$ws1.Range("F21").Formula = "=C19+1"

# Selection matches the authored workbook
$ws1.Range("G18").Select() | Out-Null
